# Delete data rows 2 through 84 (83 rows) from Sheet1, shifting remaining
# rows up. These correspond to the first 83 barcode entries ("0483"
# through "0852") that were removed from the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = $ws.Range("A2:A84").EntireRow
$rowsToDelete.Delete()

# Move/select the active cell to A2, matching the post-edit selection.
$ws.Range("A2").Select()
